# Updated cryptos list on Thu Aug 31 17:00:56 UTC 2023 with GitHub Actions
#
# The Price (column D) and Volume(1h) (column E) cells are stored as plain
# text in this sheet (e.g. "0.06510", "26.581.59") rather than numbers, so
# force the target cells to a Text number format before writing the new
# values. Otherwise Excel would auto-convert numeric-looking strings (and
# drop significant trailing/grouping characters, e.g. "0.06510" -> 0.0651).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceCells = @(
    "D2","D3","D4","D5","D6","D8","D9","D10","D11","D12","D13","D14","D15",
    "D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27",
    "D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39",
    "D40","D41","D42","D43","D44","D45","D46","D47","D48","D50","D51"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value  = "26.581.59"
$ws.Range("E2").Value  = "  -2.42%  "

$ws.Range("D3").Value  = "1.674.58"
$ws.Range("E3").Value  = "  -1.74%  "

$ws.Range("D4").Value  = "1.009"
$ws.Range("E4").Value  = "  +0.60%  "

$ws.Range("D5").Value  = "220.50"
$ws.Range("E5").Value  = "  -1.10%  "

$ws.Range("D6").Value  = "0.5199"
$ws.Range("E6").Value  = "  -1.85%  "

$ws.Range("E7").Value  = "  +0.53%  "

$ws.Range("D8").Value  = "0.06510"
$ws.Range("E8").Value  = "  -1.00%  "

$ws.Range("D9").Value  = "0.2581"
$ws.Range("E9").Value  = "  -2.58%  "

$ws.Range("D10").Value = "19.99"
$ws.Range("E10").Value = "  -3.72%  "

$ws.Range("D11").Value = "0.07693"
$ws.Range("E11").Value = "  +0.60%  "

$ws.Range("D12").Value = "1.910.65"
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("D13").Value = "1.680.85"
$ws.Range("E13").Value = "  -1.74%  "

$ws.Range("D14").Value = "4.344"
$ws.Range("E14").Value = "  -5.01%  "

$ws.Range("D15").Value = "0.5613"
$ws.Range("E15").Value = "  -1.88%  "

$ws.Range("D16").Value = "0.0₅8040"
$ws.Range("E16").Value = "  -1.61%  "

$ws.Range("D17").Value = "65.46"
$ws.Range("E17").Value = "  -2.92%  "

$ws.Range("D18").Value = "26.648.46"
$ws.Range("E18").Value = "  -2.08%  "

$ws.Range("D19").Value = "213.58"
$ws.Range("E19").Value = "  -1.18%  "

$ws.Range("D20").Value = "1.009"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("D21").Value = "4.492"
$ws.Range("E21").Value = "  -3.55%  "

$ws.Range("D22").Value = "10.15"
$ws.Range("E22").Value = "  -2.67%  "

$ws.Range("D23").Value = "5.923"
$ws.Range("E23").Value = "  -0.64%  "

$ws.Range("D24").Value = "1.010"
$ws.Range("E24").Value = "  +0.60%  "

$ws.Range("D25").Value = "143.33"
$ws.Range("E25").Value = "  +0.87%  "

$ws.Range("D26").Value = "1.729"
$ws.Range("E26").Value = "  -0.67%  "

$ws.Range("D27").Value = "0.1171"
$ws.Range("E27").Value = "  -3.87%  "

$ws.Range("D28").Value = "7.052"
$ws.Range("E28").Value = "  -2.68%  "

$ws.Range("D29").Value = "15.83"
$ws.Range("E29").Value = "  -2.66%  "

$ws.Range("D30").Value = "0.05258"
$ws.Range("E30").Value = "  -1.93%  "

$ws.Range("D31").Value = "1.272"
$ws.Range("E31").Value = "  -1.33%  "

$ws.Range("D32").Value = "3.368"
$ws.Range("E32").Value = "  -3.88%  "

$ws.Range("D33").Value = "3.254"
$ws.Range("E33").Value = "  -4.50%  "

$ws.Range("D34").Value = "1.596"
$ws.Range("E34").Value = "  -2.18%  "

$ws.Range("D35").Value = "2.772"
$ws.Range("E35").Value = "  -3.52%  "

$ws.Range("D36").Value = "2.392"
$ws.Range("E36").Value = "  -1.19%  "

$ws.Range("D37").Value = "0.9321"
$ws.Range("E37").Value = "  -1.31%  "

$ws.Range("D38").Value = "0.5751"
$ws.Range("E38").Value = "  -1.64%  "

$ws.Range("D39").Value = "1.166.27"
$ws.Range("E39").Value = "  +12.33%  "

$ws.Range("D40").Value = "0.01610"
$ws.Range("E40").Value = "  -1.16%  "

$ws.Range("D41").Value = "1.009"
$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("D42").Value = "5.714"
$ws.Range("E42").Value = "  -2.47%  "

$ws.Range("D43").Value = "0.8317"
$ws.Range("E43").Value = "  -0.74%  "

$ws.Range("D44").Value = "99.95"
$ws.Range("E44").Value = "  -0.97%  "

$ws.Range("D45").Value = "1.819.62"
$ws.Range("E45").Value = "  -1.44%  "

$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  -1.63%  "

# Rows 47 and 48 swap: Aave and Mantle trade places in the ranking.
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.4498"
$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "56.03"
$ws.Range("E48").Value = "  -3.29%  "

$ws.Range("E49").Value = "  +0.47%  "

$ws.Range("D50").Value = "7.974"
$ws.Range("E50").Value = "  -0.95%  "

$ws.Range("D51").Value = "0.05170"
$ws.Range("E51").Value = "  -1.27%  "
